# DN node 9. RES updated
#
# The RES (renewable energy source) installed-capacity data feeding the
# "Pg, Winter, S1/S2/S3" and "Pg, Summer, S1/S2/S3" sheets was refreshed,
# so the workbook needs a full recalculation. Every value cell on those
# six sheets is driven by a volatile formula of the form
#   VLOOKUP(...)*(AVERAGE('[1]Profiles, RES, <Season>'!...)*(RANDBETWEEN(95,105)/100))
# A full recalculation re-draws the RANDBETWEEN(95,105) term for each of
# those cells, producing new cached <v> results without touching any
# formula text, styles, or sheet structure.

$wb = $excel.ActiveWorkbook

# Make sure we are in automatic calculation mode, then force Excel to
# recompute every formula in the workbook (including volatile ones),
# exactly like pressing Ctrl+Alt+F9 / Ctrl+Shift+F9 would.
$excel.Calculation = -4105   # xlCalculationAutomatic
$excel.CalculateFullRebuild()
$excel.CalculateFull()
$wb.RefreshAll()
$excel.Calculate()
